# Revert "Cleaned up code. Removed DIP getAll fcn as unstable"
# -> restore the original D9 value and leave the cursor where the
#    author last left it before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore the pre-cleanup value in D9 (All / Q14 column).
$ws.Range("D9").Value = 673

# Author's selection at save time.
$ws.Range("L17").Select() | Out-Null
